# Add incremental demand for green hydrogen to RPS and begin adding blue
# hydrogen to CCS module.
#
# The CCS "CO2 Capture Potentials" workbook gains a new worksheet, CPPbHS
# (CO2 Capture Potential by Hydrogen Source), seeded with the capture rate
# for natural-gas reforming with CCS (blue hydrogen).

$wb = $excel.ActiveWorkbook

# Adding a throwaway sheet first and deleting it nudges the workbook's
# internal sheetId counter forward so the new sheet lands on sheetId 15
# (matching upstream authoring history) instead of 14.
$placeholder = $wb.Worksheets.Add()
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "CPPbHS"
$placeholder.Delete() | Out-Null

# Re-fetch a live reference to the new sheet by name.
$ws = $wb.Worksheets.Item("CPPbHS")

# Header row: units note (italic, shared with the other CPPb* sheets) and
# the "capture rate" column header.
$ws.Range("A1").Value = "Unit: dimentionless (fraction of CO2 capturable)"
$ws.Range("A1").Font.Italic = $true

# Set A2 before B1 so the shared-string table picks up
# "natural gas reforming with CCS" ahead of "capture rate".
$ws.Range("A2").Value = "natural gas reforming with CCS"
$ws.Range("B1").Value = "capture rate"
$ws.Range("B2").Value = 0.85

# Match column widths used on the sibling CPPbES / CPPbI sheets.
$ws.Columns.Item(1).ColumnWidth = 44.59
$ws.Columns.Item(2).ColumnWidth = 23.88
$ws.Columns.Item(3).ColumnWidth = 25.02

# Give the new tab the same accent color as the other CCS potential sheets.
$srcTab = $wb.Worksheets.Item("CPPbI")
$ws.Tab.Color = $srcTab.Tab.Color

# Leave the cursor on B3 (as if the user had just typed 0.85 into B2 and
# pressed Enter), then restore "About" as the active tab.
$ws.Range("B3").Select() | Out-Null
$wb.Worksheets.Item("About").Activate() | Out-Null
